# add "Save" column (H) to the s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy the header formatting from G1 (bold, centered, bordered)
# so the new "Save" header matches the existing header styling, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Data cells: plain numeric values, no special style (matches F/G columns' data rows)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
